# Update benchmark: 2026-01-07 06:43:19 UTC
# Moves several per-bank fee figures from one bank column to another
# (columns C/F/G/H across several rows), tweaks a couple of values, and
# fills in a few previously-empty cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BENCHMARK")

# --- Row 3: EFT - HESAPTAN EFT - Şube ---
$ws.Range("C3").ClearContents()
$ws.Range("F3").Value = "30,46 TL - 60,94 TL - 609,43 TL"
$ws.Range("G3").Value = "30,46 TL - 60,94 TL - 609,43 TL"
$ws.Range("H3").ClearContents()

# --- Row 4: EFT - HESAPTAN EFT - ATM ---
$ws.Range("C4").ClearContents()
$ws.Range("F4").Value = "30,46 TL - 60,94 TL - 609,43 TL"
$ws.Range("G4").Value = "30,46 TL - 60,94 TL - 609,43 TL"
$ws.Range("H4").ClearContents()

# --- Row 5: EFT - HESAPTAN EFT - Mobil ---
$ws.Range("C5").ClearContents()
$ws.Range("F5").Value = "30,46 TL - 60,94 TL - 609,43 TL"
$ws.Range("G5").Value = "30,46 TL - 60,94 TL - 609,43 TL"
$ws.Range("H5").ClearContents()

# --- Row 6: DÜZENLİ EFT ---
$ws.Range("C6").ClearContents()
$ws.Range("G6").Value = "6,09 TL - 12,19 TL - 152,35 TL"
$ws.Range("H6").ClearContents()

# --- Row 8: HAVALE - HESAPTAN HAVALE - Şube ---
$ws.Range("C8").ClearContents()
$ws.Range("F8").Value = "15,23 TL - 30,47 TL - 304,71 TL"
$ws.Range("G8").Value = "15,23 TL - 30,47 TL - 211,05 TL"
$ws.Range("H8").ClearContents()

# --- Row 9: HAVALE - HESAPTAN HAVALE - ATM ---
$ws.Range("C9").ClearContents()
$ws.Range("F9").Value = "15,23 TL - 30,47 TL - 304,71 TL"
$ws.Range("G9").Value = "15,23 TL - 30,47 TL - 211,05 TL"
$ws.Range("H9").ClearContents()

# --- Row 10: HAVALE - HESAPTAN HAVALE - Mobil ---
$ws.Range("C10").ClearContents()
$ws.Range("F10").Value = "15,23 TL - 30,47 TL - 304,71 TL"
$ws.Range("G10").Value = "15,23 TL - 30,47 TL - 211,05 TL"
$ws.Range("H10").ClearContents()

# --- Row 11: DÜZENLİ HAVALE ---
$ws.Range("C11").ClearContents()
$ws.Range("G11").Value = "3,04 TL - 6,09 TL - 76,17 TL"
$ws.Range("H11").ClearContents()

# --- Row 12: SWIFT - GİDEN SWIFT ---
$ws.Range("C12").ClearContents()
$ws.Range("G12").Value = "Şube (Kasadan): %0,5; Şube (Hesaptan): %0,75; İnternet: 15 USD"

# --- Row 13: GELEN SWIFT ---
$ws.Range("C13").ClearContents()
$ws.Range("E13").Value = "Hesaba: Asgari 1 TL | Azami 8.700 TL"
$ws.Range("F13").Value = "Hesaba: Asgari 300 TL | Azami 3.080 TL"
$ws.Range("H13").ClearContents()

# --- Row 14: GİDEN SWIFT - Mobil ---
$ws.Range("C14").ClearContents()
$ws.Range("F14").Value = "1.952,38 TL - 9.523,81 TL"
$ws.Range("G14").Value = "6.300 TL - 6,09 TL"
$ws.Range("H14").ClearContents()

# --- Row 24: SENET - SENET TAHSİLE ALMA ---
$ws.Range("D24").Value = "476,2 TL"
$ws.Range("J24").Value = "375 TL"

# --- Row 25: MUAMELESİZ SENET İADESİ ---
$ws.Range("D25").Value = "428,58 TL"
$ws.Range("J25").Value = "375 TL"
